# "Feature List.xlsx" update — Aug revision
# Restructures the Sprint plan: inserts a new "Sprint 3" row (shifting the
# remaining Sprint rows down by one, through a new "Sprint 14" at the end),
# splits the old combined "sign in / upload" and "map" requirement
# descriptions into separate, more granular user stories, and adds a new
# "Application - route lnstructions" column with its own set of sprint
# deliverables (beacon/QR/barcode/NFC detection, Organisation code
# download, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "route lnstructions" header column -------------------------------
$ws.Range("F1").Value = "Application - route lnstructions"

# --- Row 2 ("What is is?") — split the old combined description ----------
# the old B2/C2/D2/E2 text is replaced by four focused stories, spread
# across rows 2-4 below and the new E-column beacon/QR/barcode/NFC rows.
$ws.Range("B2").Value = "The organisation can sign in"
$ws.Range("C2:E2").ClearContents()

# --- Row 3 ("Sprint 1") ----------------------------------------------------
$ws.Range("B3").Value = "Organisation can upload their maps"
$ws.Range("C3").Value = "Read in the Organisations Code to download data to the application"

# --- Row 4 ("Sprint 2") ----------------------------------------------------
$ws.Range("B4").Value = "Organisation can upload information for navigation to work"
$ws.Range("C4").Value = "To show the map after download of data and option to select where the user would like to go."

# --- Row 5 — previously blank, now "Sprint 3" ------------------------------
$ws.Range("A5").Value = "Sprint 3"
$ws.Range("D5").Value = "Marker to display where you are"
$ws.Range("F5").Value = "Show list of stepts to get from where you are to where you are going"

# --- Rows 6-9: Sprint numbers shift down by one & new beacon/QR/barcode/NFC detail --
$ws.Range("A6").Value = "Sprint 4"
$ws.Range("E6").Value = "Detect where you are by using bluetooth beacons"

$ws.Range("A7").Value = "Sprint 5"
$ws.Range("E7").Value = "Detect where you are by using NFC"

$ws.Range("A8").Value = "Sprint 6"
$ws.Range("E8").Value = "Detect where you are by using QR Codes"

$ws.Range("A9").Value = "Sprint 7"
$ws.Range("E9").Value = "Detect where you are by using Bar Codes"

# --- Rows 10-16: remaining Sprint numbers shift down by one, plus a new Sprint 14 --
$ws.Range("A10").Value = "Sprint 8"
$ws.Range("A11").Value = "Sprint 9"
$ws.Range("A12").Value = "Sprint 10"
$ws.Range("A13").Value = "Sprint 11"
$ws.Range("A14").Value = "Sprint 12"
$ws.Range("A15").Value = "Sprint 13"
$ws.Range("A16").Value = "Sprint 14"

# --- Carry over the cell formatting (fill/border/alignment) that already
# exists elsewhere on the sheet onto the newly-populated cells, matching
# the look of their neighbours -------------------------------------------
$ws.Range("D3").Copy()
$ws.Range("C2:H2").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("F5").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("E6:E9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Leave the active selection where the author left it ------------------
$ws.Range("E9").Select() | Out-Null
